# Updates cryptos list figures (price + 1h volume change) plus a row swap
# between "Maker" and "EnergySwap", matching the Jan 29 2024 GitHub Actions
# data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rangeRef, $val) {
    # Force the cell to stay a text value even when the string looks like a
    # number (e.g. "1.00", "35.51"), then drop back to the default style so
    # we don't leave a stray NumberFormat behind on the cell.
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# row -> (Price, Volume(1h)) ; $null means "leave Price unchanged"
$rows = @(
    @{ Row = 2;  D = "42.447.39";  E = "  +0.78%  " },
    @{ Row = 3;  D = "2.276.84";   E = "  -0.13%  " },
    @{ Row = 4;  D = "1.00";       E = "  -0.08%  " },
    @{ Row = 5;  D = "306.62";     E = "  +0.45%  " },
    @{ Row = 6;  D = "97.97";      E = "  +3.29%  " },
    @{ Row = 7;  D = "0.529";      E = "  -0.55%  " },
    @{ Row = 8;  D = $null;        E = "  -0.03%  " },
    @{ Row = 9;  D = "0.497";      E = "  +0.78%  " },
    @{ Row = 10; D = "35.51";      E = "  -0.40%  " },
    @{ Row = 11; D = $null;        E = "  -1.31%  " },
    @{ Row = 12; D = $null;        E = "  +0.01%  " },
    @{ Row = 13; D = $null;        E = "  +3.11%  " },
    @{ Row = 14; D = "2.629.48";   E = "  -0.15%  " },
    @{ Row = 15; D = $null;        E = "  +2.82%  " },
    @{ Row = 16; D = "2.276.85";   E = "  +0.09%  " },
    @{ Row = 17; D = "0.798";      E = "  +0.34%  " },
    @{ Row = 18; D = "42.282.70";  E = "  +0.57%  " },
    @{ Row = 19; D = "12.55";      E = "  -2.12%  " },
    @{ Row = 20; D = "0.0₃0909";   E = "  -0.95%  " },
    @{ Row = 21; D = $null;        E = "  +0.67%  " },
    @{ Row = 22; D = "68.30";      E = "  +0.28%  " },
    @{ Row = 23; D = "239.44";     E = "  -1.79%  " },
    @{ Row = 24; D = "2.59";       E = "  -0.59%  " },
    @{ Row = 25; D = "1.97";       E = "  +1.00%  " },
    @{ Row = 26; D = $null;        E = "  +0.13%  " },
    @{ Row = 27; D = "23.73";      E = "  -1.67%  " },
    @{ Row = 28; D = "38.34";      E = "  +5.16%  " },
    @{ Row = 29; D = $null;        E = "  -1.67%  " },
    @{ Row = 30; D = "2.12";       E = "  +0.86%  " },
    @{ Row = 31; D = "161.89";     E = "  +0.16%  " },
    @{ Row = 32; D = "5.27";       E = "  -1.67%  " },
    @{ Row = 34; D = "3.20";       E = "  +3.79%  " },
    @{ Row = 35; D = "0.0741";     E = "  -1.65%  " },
    @{ Row = 36; D = "17.60";      E = "  +2.72%  " },
    @{ Row = 37; D = $null;        E = "  -0.48%  " },
    @{ Row = 38; D = $null;        E = "  -3.13%  " },
    @{ Row = 39; D = $null;        E = "  +1.10%  " },
    @{ Row = 40; D = $null;        E = "  -1.49%  " },
    @{ Row = 41; D = "4.11";       E = "  -1.86%  " },
    @{ Row = 42; D = $null;        E = "  +2.39%  " },
    @{ Row = 45; D = "0.0283";     E = "  -0.35%  " },
    @{ Row = 46; D = $null;        E = "  -2.36%  " },
    @{ Row = 47; D = $null;        E = "  -2.64%  " },
    @{ Row = 48; D = "53.85";      E = "  +0.47%  " },
    @{ Row = 49; D = "92.75";      E = "  +0.54%  " },
    @{ Row = 50; D = "72.23";      E = "  -0.21%  " },
    @{ Row = 51; D = $null;        E = "  -1.24%  " }
)

foreach ($item in $rows) {
    $r = $item.Row
    if ($null -ne $item.D) {
        Set-TextCell "D$r" $item.D
    }
    $ws.Range("E$r").Value = $item.E
}

# Rows 43/44 swap: "Maker" and "EnergySwap" trade places (with refreshed
# price / volume figures) rather than simply updating in place.
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D43" "19.08"
$ws.Range("E43").Value = "  -2.68%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D44" "1.950.82"
$ws.Range("E44").Value = "  -3.39%  "
